# Add a new worksheet "ʤ" (U+02A4) after the "data" sheet, populate it with
# a small a/b/c header row + 1/2/3 data row, and restore "data" as the
# active/selected sheet so the workbook's active-tab bookkeeping is
# unchanged from before the edit.

$wb = $excel.ActiveWorkbook

$dataSheet = $wb.Worksheets.Item("data")

# Insert the new sheet directly after "data" (becomes the last tab).
$ws = $wb.Worksheets.Add($null, $dataSheet)
$ws.Name = [char]0x02A4

$ws.Range("A1").Value = "a"
$ws.Range("B1").Value = "b"
$ws.Range("C1").Value = "c"
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 3

# Adding the sheet makes it active; put the selection back on "data".
$dataSheet.Activate()
